# Update Lotofacil_edt2.xlsx: append the latest 5 lottery draws (concursos
# 3488-3492) to the bottom of the results table, and refresh the sheet view
# (freeze the header row, scroll down, select the cell next to the newest
# row) to mirror how the spreadsheet was left after the update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# New lottery draws to append (Concurso, Bola1..Bola15)
$newRows = @(
    @(3488,1,3,4,5,6,8,9,11,12,13,14,15,17,22,25),
    @(3489,1,2,5,8,9,11,14,16,17,20,21,22,23,24,25),
    @(3490,2,3,4,7,8,11,13,14,15,16,18,19,21,23,25),
    @(3491,1,2,4,8,9,10,12,13,15,17,21,22,23,24,25),
    @(3492,2,3,4,8,9,10,13,17,18,19,20,21,22,23,25)
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

# Refresh the window layout: freeze the header row, scroll the view near
# the bottom of the data, and leave the selection on F420.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 395
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("F420").Select()
